# Updates the cryptocurrency price/volume table on Sheet1 to the latest
# scraped snapshot. Most D-column (Price) cells hold numeric-looking text
# (e.g. "0.999", "39.41") that must stay text, exactly as scraped - so we
# force NumberFormat "@" (Text) on those specific cells before writing the
# value, which keeps Excel's COM layer from auto-coercing them to numbers.
# Three rows (37/38, 40/41/42, 47/48) had their coin entries re-ranked, so
# their B (Coin), C (Link), D (Price) and E (Volume) cells are fully
# rewritten rather than just the numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.957.35"
$ws.Range("E2").Value = "  +8.51%  "
$ws.Range("D3").Value = "3.217.85"
$ws.Range("E3").Value = "  +3.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "398.20"
$ws.Range("E5").Value = "  +3.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.16"
$ws.Range("E6").Value = "  +6.31%  "
$ws.Range("E7").Value = "  +3.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.623"
$ws.Range("E9").Value = "  +6.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.41"
$ws.Range("E10").Value = "  +6.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0908"
$ws.Range("E11").Value = "  +6.20%  "
$ws.Range("E12").Value = "  +2.16%  "
$ws.Range("D13").Value = "3.723.93"
$ws.Range("E13").Value = "  +3.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.12"
$ws.Range("E14").Value = "  +3.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.13"
$ws.Range("E15").Value = "  +3.03%  "
$ws.Range("D16").Value = "3.209.61"
$ws.Range("E16").Value = "  +3.51%  "
$ws.Range("E17").Value = "  +5.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.73"
$ws.Range("E18").Value = "  -1.55%  "
$ws.Range("D19").Value = "55.736.72"
$ws.Range("E19").Value = "  +8.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.37"
$ws.Range("E20").Value = "  +2.91%  "
$ws.Range("E21").Value = "  +7.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.12"
$ws.Range("E22").Value = "  +5.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "304.84"
$ws.Range("E23").Value = "  +14.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.30"
$ws.Range("E24").Value = "  +7.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.25"
$ws.Range("E25").Value = "  +2.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.22"
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "28.28"
$ws.Range("E27").Value = "  +4.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.52"
$ws.Range("E28").Value = "  +5.36%  "
$ws.Range("E29").Value = "  +4.63%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("E31").Value = "  +4.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.27"
$ws.Range("E32").Value = "  +8.55%  "
$ws.Range("E33").Value = "  +2.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "36.21"
$ws.Range("E34").Value = "  +2.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.12"
$ws.Range("E35").Value = "  +2.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.40"
$ws.Range("E36").Value = "  +2.67%  "
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.11"
$ws.Range("E37").Value = "  +23.19%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.49"
$ws.Range("E39").Value = "  +4.22%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.93"
$ws.Range("E40").Value = "  +2.69%  "
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.04"
$ws.Range("E41").Value = "  +7.19%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "133.45"
$ws.Range("E42").Value = "  +3.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.20"
$ws.Range("E43").Value = "  +3.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.287"
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("E45").Value = "  +2.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.41"
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.49"
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.11"
$ws.Range("E48").Value = "  +45.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.09"
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("D50").Value = "2.137.94"
$ws.Range("E50").Value = "  +3.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0366"
$ws.Range("E51").Value = "  +13.20%  "
